# Auto-generated edit script applying the cryptos.xlsx price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.515.70"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.624.47"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.52"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").Value = "3.100.01"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +14.16%  "
$ws.Range("D15").Value = "60.513.76"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "2.626.90"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.71"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.03%  "
$ws.Range("E28").Value = "  +13.27%  "
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +12.01%  "
$ws.Range("E36").Value = "  +8.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.65"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "331.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +12.53%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  +5.46%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.59%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0996"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.610"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0245"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.84%  "
